$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new data row at row 7 (Compleat Sculptor order), inheriting the
#     fully-formatted look (date / wrap-text / currency styles) of row 6 ---
$ws.Rows("6:6").Copy()
$ws.Rows("7:7").Insert()

# --- The sheet keeps 36 total rows: one of the two still-unused placeholder
#     rows (the one that was row 9) is removed so the row count balances out ---
$ws.Rows("10:10").Delete()

# --- Populate the newly inserted row 7 ---
$ws.Range("A7").Value = "2/20/2018"
$ws.Range("B7").Value = "Compleat Sculptor"
$ws.Range("C7").Value = "XTC Smooth-On Epoxy Die"
$ws.Range("D7").Value = "The Compleat Sculptor"
$ws.Range("E7").Value = 25

# --- Populate row 10 (Waterjet order), matching the date formatting used by
#     the rest of the table's date column ---
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = "4/12/2018"
$ws.Range("B10").Value = "Waterjet Order"
$ws.Range("C10").Value = "5 Sealing Discs Waterjet Order"
$ws.Range("D10").Value = "RPI MILL"
$ws.Range("E10").Value = 5

# --- Widen the Supplier(s) column to fit the new, longer supplier names ---
$ws.Columns("D:D").ColumnWidth = 24.86

# --- Match the author's final cursor position ---
$ws.Range("E11").Select()
